$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw data values for "Bohemian Rhapsody" (row 2)
$ws.Range("B2").Value = 2535
$ws.Range("C2").Value = 248

# Force recalculation of all dependent formulas (sums, regression, correlation, chart caches)
$excel.CalculateFull()

# Move the active selection to E4, as in the edited workbook
$ws.Range("E4").Select()
